# Analysis.xlsx edit: new Max-Flow implementation data added / re-run,
# per commit message "Added some new max flow implementations, for
# further testing and analysis".
#
# Updates the raw trial timings (columns D:H) for a few algorithm rows;
# the AVERAGE()/LOG() formulas in the same rows recompute automatically.
# Also repositions the two result charts and refreshes the active
# selection / zoom level on the sheet, matching the author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated raw timing values (D:H) -----------------------------------
# Row 5  (n = 10..17 input)
$ws.Range("D5").Value = 918400
$ws.Range("E5").Value = 827700
$ws.Range("F5").Value = 826900
$ws.Range("G5").Value = 826900

# Row 7  (n = 34..65 input)
$ws.Range("E7").Value = 910700
$ws.Range("H7").Value = 1011000

# Row 10 (n = 258..513 input)
$ws.Range("E10").Value = 6486700

# Row 14 (n = 6..9 input, second series)
$ws.Range("D14").Value = 752500
$ws.Range("E14").Value = 688000
$ws.Range("F14").Value = 717900
$ws.Range("G14").Value = 774600

# --- Reposition the two charts ------------------------------------------
$co1 = $ws.ChartObjects(1)
$co1.Left = 1231.1158313852505
$co1.Top = 252.25708661417322
$co1.Width = 534.9421238312086
$co1.Height = 392.51590551181105

$co2 = $ws.ChartObjects(2)
$co2.Left = 142.15638102854712
$co2.Top = 384.0658267716535
$co2.Width = 442.84723763530404
$co2.Height = 400.3514173228346

# --- Refresh active selection / zoom -------------------------------------
$ws.Range("G14").Select()
$excel.ActiveWindow.Zoom = 115
